# Update the spreadsheet per the diff:
# - delete the entire row 3
# - delete columns BC and BD (they only had data in rows 1-3, which are now gone/removed)
# - update a handful of odds values in row 2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in row 2
$ws.Range("I2").Value = 4.5
$ws.Range("Q2").Value = 1.7
$ws.Range("R2").Value = 2.1
$ws.Range("AD2").Value = 7
$ws.Range("AY2").Value = 29

# Delete entire row 3
$ws.Rows.Item(3).Delete()

# Delete columns BC and BD entirely
$ws.Range("BC:BD").Delete()
